$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (AVG_TIME_TO_HF)
$ws.Range("B2").Value = 0.4854213574091631
$ws.Range("C2").Value = 0.9903342299228689
$ws.Range("D2").Value = 0.567746984583792
$ws.Range("G2").Value = 0.4821145882335259
$ws.Range("H2").Value = 0.992

# Row 3 (AVG_TIME_TO_MI)
$ws.Range("B3").Value = 0.2403775317171263
$ws.Range("C3").Value = 0.9953018150975473
$ws.Range("D3").Value = 0.3872771956938404
$ws.Range("G3").Value = 0.4821145882335259
$ws.Range("H3").Value = 0.992

# Row 4 (AVG_TIME_TO_ANGINA)
$ws.Range("B4").Value = 0.3477187780492769
$ws.Range("C4").Value = 0.9933106639122377
$ws.Range("D4").Value = 0.4697298980578373
$ws.Range("G4").Value = 0.4821145882335259
$ws.Range("H4").Value = 0.992

# Row 5 (AVG_TIME_TO_STROKE)
$ws.Range("B5").Value = 0.383890278696353
$ws.Range("C5").Value = 0.9924302092967611
$ws.Range("D5").Value = 0.4868141901908184
$ws.Range("G5").Value = 0.4821145882335259
$ws.Range("H5").Value = 0.992
